# Refresh the cryptos price/volume snapshot (Price = column D, Volume(1h) = column E).
# Values are written as literal text (matching the sheet's existing inline-string
# cells); a leading apostrophe is used for prices that would otherwise be
# auto-parsed by Excel as numbers, so they stay plain text like the originals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.488.61"
$ws.Range("E2").Value = "  -4.22%  "
$ws.Range("D3").Value = "3.107.97"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("D5").Value = "'549.94"
$ws.Range("E5").Value = "  -4.52%  "
$ws.Range("D6").Value = "'137.85"
$ws.Range("E6").Value = "  -10.81%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.103.68"
$ws.Range("E8").Value = "  -4.03%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "  -10.03%  "
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").Value = "'35.64"
$ws.Range("E13").Value = "  -5.60%  "
$ws.Range("E14").Value = "  -7.61%  "
$ws.Range("D15").Value = "3.609.52"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "63.475.80"
$ws.Range("E16").Value = "  -4.34%  "
$ws.Range("D17").Value = "'0.112"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("D18").Value = "3.104.46"
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("E19").Value = "  -5.18%  "
$ws.Range("D20").Value = "'492.35"
$ws.Range("E20").Value = "  -12.47%  "
$ws.Range("E21").Value = "  -5.12%  "
$ws.Range("D22").Value = "'0.721"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "'7.26"
$ws.Range("E23").Value = "  -7.82%  "
$ws.Range("D24").Value = "'79.28"
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D25").Value = "'12.42"
$ws.Range("E25").Value = "  -8.25%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'8.50"
$ws.Range("E27").Value = "  -9.22%  "
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("E29").Value = "  -11.85%  "
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").Value = "'26.71"
$ws.Range("E31").Value = "  -4.54%  "
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("E33").Value = "  -9.01%  "
$ws.Range("D34").Value = "'59.13"
$ws.Range("E34").Value = "  +6.69%  "
$ws.Range("D35").Value = "'519.39"
$ws.Range("E35").Value = "  -8.36%  "
$ws.Range("D36").Value = "'6.02"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("E37").Value = "  -10.34%  "
$ws.Range("D38").Value = "'0.0407"
$ws.Range("E38").Value = "  -10.57%  "
$ws.Range("D39").Value = "3.159.02"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'0.0807"
$ws.Range("E40").Value = "  -7.03%  "
$ws.Range("D41").Value = "'0.120"
$ws.Range("E41").Value = "  -5.54%  "
$ws.Range("D42").Value = "'8.20"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("E43").Value = "  -12.52%  "
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("E46").Value = "  -10.49%  "
$ws.Range("D47").Value = "'25.13"
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("D48").Value = "'121.31"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "0.0₃0509"
$ws.Range("E50").Value = "  -9.20%  "
$ws.Range("E51").Value = "  -9.74%  "
